$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values per diff
$ws.Range("B2").Value = 5.5
$ws.Range("C2").Value = 11
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 9.5
$ws.Range("C4").Value = 1.25

# Update active selection to C5
$ws.Range("C5").Select()
